# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 13:22"

# Row 14 (Suiza) - updated case counts
$ws.Range("B14").Value = 24900
$ws.Range("C14").Value = 349
$ws.Range("E14").Value = 12797

# Row 19 (Austria) - updated case counts
$ws.Range("B19").Value = 13744
$ws.Range("C19").Value = 184
$ws.Range("E19").Value = 6803

# Rows 47-49: Catar overtakes Republica Dominicana and Tailandia in ranking
# (data refresh re-sorted by total cases, so the country names shift down
#  a row while Catar gets fresh totals)
$ws.Range("A47").Value = "Catar"
$ws.Range("B47").Value = 2728
$ws.Range("C47").Value = 216
$ws.Range("D47").Value = 247
$ws.Range("E47").Value = 2475
$ws.Range("F47").Value = 37
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 6

$ws.Range("A48").Value = "Republica Dominicana"
$ws.Range("B48").Value = 2620
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 98
$ws.Range("E48").Value = 2396
$ws.Range("F48").Value = 147
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 126

$ws.Range("A49").Value = "Tailandia"
$ws.Range("B49").Value = 2518
$ws.Range("C49").Value = 45
$ws.Range("D49").Value = 1135
$ws.Range("E49").Value = 1348
$ws.Range("F49").Value = 61
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 6

# Row 73 (Armenia) - updated case counts
$ws.Range("B73").Value = 967
$ws.Range("C73").Value = 30
$ws.Range("E73").Value = 781

# Row 112 (Georgia) - updated case counts
$ws.Range("D112").Value = 56
$ws.Range("E112").Value = 175

# Row 211 (Islas Virgenes Britanicas) - updated case counts
$ws.Range("D211").Value = 2
$ws.Range("E211").Value = 1

# Rows 212-213: Bonaire, San Eustaquio y Saba swaps position with
# Papua Nueva Guinea in the country ordering
$ws.Range("A212").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A213").Value = "Papua Nueva Guinea"
